# Auto-generated Excel COM-interop script to apply the Línea 141 - 628 schedule update
$wb = $excel.ActiveWorkbook

# ----- Sheet 1: LP1912 -----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 05:26:08"
$ws1.Range("A3").Value = "Total filas: 36"

$ws1.Cells.Item(14,1).Value = "05:26:08"
$ws1.Cells.Item(14,2).Value = "05:28"
$ws1.Cells.Item(14,3).Value = "14_ABASTO"
$ws1.Cells.Item(14,4).Value = 2
$ws1.Cells.Item(14,5).Value = "LP1912"

$ws1.Cells.Item(15,1).Value = "04:18:02"
$ws1.Cells.Item(15,2).Value = "05:34"
$ws1.Cells.Item(15,3).Value = "14_ABASTO"
$ws1.Cells.Item(15,4).Value = 76
$ws1.Cells.Item(15,5).Value = "LP1912"

$ws1.Cells.Item(16,1).Value = "03:45:25"
$ws1.Cells.Item(16,2).Value = "05:34"
$ws1.Cells.Item(16,3).Value = "215B_EL PATO"
$ws1.Cells.Item(16,4).Value = 109
$ws1.Cells.Item(16,5).Value = "LP1912"

$ws1.Cells.Item(17,1).Value = "04:18:02"
$ws1.Cells.Item(17,2).Value = "05:35"
$ws1.Cells.Item(17,3).Value = "215B_EL PATO"
$ws1.Cells.Item(17,4).Value = 77
$ws1.Cells.Item(17,5).Value = "LP1912"

$ws1.Cells.Item(18,1).Value = "03:45:25"
$ws1.Cells.Item(18,2).Value = "05:37"
$ws1.Cells.Item(18,3).Value = "14_ABASTO"
$ws1.Cells.Item(18,4).Value = 112
$ws1.Cells.Item(18,5).Value = "LP1912"

$ws1.Cells.Item(19,1).Value = "04:18:02"
$ws1.Cells.Item(19,2).Value = "05:46"
$ws1.Cells.Item(19,3).Value = "15_ABASTO"
$ws1.Cells.Item(19,4).Value = 88
$ws1.Cells.Item(19,5).Value = "LP1912"

$ws1.Cells.Item(20,1).Value = "04:45:05"
$ws1.Cells.Item(20,2).Value = "06:04"
$ws1.Cells.Item(20,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(20,4).Value = 79
$ws1.Cells.Item(20,5).Value = "LP1912"

$ws1.Cells.Item(21,1).Value = "04:18:02"
$ws1.Cells.Item(21,2).Value = "06:05"
$ws1.Cells.Item(21,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(21,4).Value = 107
$ws1.Cells.Item(21,5).Value = "LP1912"

$ws1.Cells.Item(22,1).Value = "04:56:49"
$ws1.Cells.Item(22,2).Value = "06:11"
$ws1.Cells.Item(22,3).Value = "215A_EL PATO"
$ws1.Cells.Item(22,4).Value = 75
$ws1.Cells.Item(22,5).Value = "LP1912"

$ws1.Cells.Item(23,1).Value = "04:18:02"
$ws1.Cells.Item(23,2).Value = "06:12"
$ws1.Cells.Item(23,3).Value = "215A_EL PATO"
$ws1.Cells.Item(23,4).Value = 114
$ws1.Cells.Item(23,5).Value = "LP1912"

$ws1.Cells.Item(24,1).Value = "04:18:02"
$ws1.Cells.Item(24,2).Value = "06:14"
$ws1.Cells.Item(24,3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(24,4).Value = 116
$ws1.Cells.Item(24,5).Value = "LP1912"

$ws1.Cells.Item(25,1).Value = "04:45:05"
$ws1.Cells.Item(25,2).Value = "06:21"
$ws1.Cells.Item(25,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(25,4).Value = 96
$ws1.Cells.Item(25,5).Value = "LP1912"

$ws1.Cells.Item(26,1).Value = "04:45:05"
$ws1.Cells.Item(26,2).Value = "06:27"
$ws1.Cells.Item(26,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(26,4).Value = 102
$ws1.Cells.Item(26,5).Value = "LP1912"

$ws1.Cells.Item(27,1).Value = "04:56:49"
$ws1.Cells.Item(27,2).Value = "06:29"
$ws1.Cells.Item(27,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(27,4).Value = 93
$ws1.Cells.Item(27,5).Value = "LP1912"

$ws1.Cells.Item(28,1).Value = "04:45:05"
$ws1.Cells.Item(28,2).Value = "06:30"
$ws1.Cells.Item(28,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(28,4).Value = 105
$ws1.Cells.Item(28,5).Value = "LP1912"

$ws1.Cells.Item(29,1).Value = "04:45:05"
$ws1.Cells.Item(29,2).Value = "06:31"
$ws1.Cells.Item(29,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(29,4).Value = 106
$ws1.Cells.Item(29,5).Value = "LP1912"

$ws1.Cells.Item(30,1).Value = "04:45:05"
$ws1.Cells.Item(30,2).Value = "06:44"
$ws1.Cells.Item(30,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(30,4).Value = 119
$ws1.Cells.Item(30,5).Value = "LP1912"

$ws1.Cells.Item(31,1).Value = "04:56:49"
$ws1.Cells.Item(31,2).Value = "06:46"
$ws1.Cells.Item(31,3).Value = "215C_EL PATO"
$ws1.Cells.Item(31,4).Value = 110
$ws1.Cells.Item(31,5).Value = "LP1912"

$ws1.Cells.Item(32,1).Value = "05:26:08"
$ws1.Cells.Item(32,2).Value = "06:47"
$ws1.Cells.Item(32,3).Value = "215C_EL PATO"
$ws1.Cells.Item(32,4).Value = 81
$ws1.Cells.Item(32,5).Value = "LP1912"

$ws1.Cells.Item(33,1).Value = "05:26:08"
$ws1.Cells.Item(33,2).Value = "07:00"
$ws1.Cells.Item(33,3).Value = "14_ABASTO"
$ws1.Cells.Item(33,4).Value = 94
$ws1.Cells.Item(33,5).Value = "LP1912"

$ws1.Cells.Item(34,1).Value = "05:26:08"
$ws1.Cells.Item(34,2).Value = "07:05"
$ws1.Cells.Item(34,3).Value = "15_ABASTO"
$ws1.Cells.Item(34,4).Value = 99
$ws1.Cells.Item(34,5).Value = "LP1912"

$ws1.Cells.Item(35,1).Value = "05:26:08"
$ws1.Cells.Item(35,2).Value = "07:05"
$ws1.Cells.Item(35,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(35,4).Value = 99
$ws1.Cells.Item(35,5).Value = "LP1912"

$ws1.Cells.Item(36,1).Value = "05:26:08"
$ws1.Cells.Item(36,2).Value = "07:06"
$ws1.Cells.Item(36,3).Value = "10_OLMOS"
$ws1.Cells.Item(36,4).Value = 100
$ws1.Cells.Item(36,5).Value = "LP1912"

$ws1.Cells.Item(37,1).Value = "05:26:08"
$ws1.Cells.Item(37,2).Value = "07:07"
$ws1.Cells.Item(37,3).Value = "225_GOMEZ"
$ws1.Cells.Item(37,4).Value = 101
$ws1.Cells.Item(37,5).Value = "LP1912"

$ws1.Cells.Item(38,1).Value = "05:26:08"
$ws1.Cells.Item(38,2).Value = "07:11"
$ws1.Cells.Item(38,3).Value = "215A_EL PATO"
$ws1.Cells.Item(38,4).Value = 105
$ws1.Cells.Item(38,5).Value = "LP1912"

$ws1.Cells.Item(39,1).Value = "05:26:08"
$ws1.Cells.Item(39,2).Value = "07:16"
$ws1.Cells.Item(39,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(39,4).Value = 110
$ws1.Cells.Item(39,5).Value = "LP1912"

$ws1.Cells.Item(40,1).Value = "05:26:08"
$ws1.Cells.Item(40,2).Value = "07:21"
$ws1.Cells.Item(40,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(40,4).Value = 115
$ws1.Cells.Item(40,5).Value = "LP1912"

$ws1.Cells.Item(41,1).Value = "05:26:08"
$ws1.Cells.Item(41,2).Value = "07:23"
$ws1.Cells.Item(41,3).Value = "10_OLMOS"
$ws1.Cells.Item(41,4).Value = 117
$ws1.Cells.Item(41,5).Value = "LP1912"

# ----- Sheet 2: LP1912-215 -----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 05:26:08"
$ws2.Range("A3").Value = "Total filas: 8"

$ws2.Cells.Item(12,1).Value = "05:26:08"
$ws2.Cells.Item(12,2).Value = "06:47"
$ws2.Cells.Item(12,3).Value = "215C_EL PATO"
$ws2.Cells.Item(12,4).Value = 81
$ws2.Cells.Item(12,5).Value = "LP1912"

$ws2.Cells.Item(13,1).Value = "05:26:08"
$ws2.Cells.Item(13,2).Value = "07:11"
$ws2.Cells.Item(13,3).Value = "215A_EL PATO"
$ws2.Cells.Item(13,4).Value = 105
$ws2.Cells.Item(13,5).Value = "LP1912"

# ----- Sheet 3: 6203-6173 -----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 05:26:08"
$ws3.Range("A3").Value = "Total filas: 7"

$ws3.Cells.Item(12,1).Value = "05:26:08"
$ws3.Cells.Item(12,2).Value = "07:00"
$ws3.Cells.Item(12,3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(12,4).Value = 94
$ws3.Cells.Item(12,5).Value = "L6173"

